$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (Last Name, First Name, Class, Equip) shifts one column to
# the left: previously B5:E5, now A5:D5.
$lastName  = $ws.Range("B5").Value2
$firstName = $ws.Range("C5").Value2
$class     = $ws.Range("D5").Value2
$equip     = $ws.Range("E5").Value2

$ws.Range("A5").Value2 = $lastName
$ws.Range("B5").Value2 = $firstName
$ws.Range("C5").Value2 = $class
$ws.Range("D5").Value2 = $equip
$ws.Range("E5").ClearContents()
